# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Map: sheet name -> hashtable of row -> new F value
$updates = @{
    "展览" = @{
        2  = 7031
        4  = 462
        7  = 151
        11 = 52
        13 = 446
        15 = 1828
        17 = 3631
        23 = 2260
        25 = 256
        32 = 260
        33 = 99
    }
    "全部类型" = @{
        2  = 7031
        4  = 462
        8  = 151
        12 = 52
        14 = 446
        16 = 1828
        18 = 3631
        24 = 2260
        26 = 256
        33 = 260
        34 = 99
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowUpdates = $updates[$sheetName]
    foreach ($row in $rowUpdates.Keys) {
        $ws.Range("F$row").Value = $rowUpdates[$row]
    }
}
